$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the coin price (D) / 1h volume change (E) columns with the
# latest scrape. D values that look like plain numbers ("178.20",
# "0.0000282", ...) get a leading apostrophe so Excel keeps them as text
# (matching the original inlineStr cells) instead of silently coercing
# them into numbers and dropping trailing zeros / using exponent notation.

$ws.Range("D2").Value = "69.227.41"
$ws.Range("E2").Value = "  +1.24%  "
$ws.Range("D3").Value = "3.416.34"
$ws.Range("E3").Value = "  +1.89%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'581.85"
$ws.Range("E5").Value = "  -0.52%  "
$ws.Range("D6").Value = "'178.20"
$ws.Range("E6").Value = "  +0.41%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'0.592"
$ws.Range("E8").Value = "  +0.28%  "
$ws.Range("D9").Value = "'0.198"
$ws.Range("E9").Value = "  +7.53%  "
$ws.Range("D10").Value = "'0.585"
$ws.Range("E10").Value = "  +0.50%  "
$ws.Range("D11").Value = "'48.25"
$ws.Range("E11").Value = "  +0.63%  "
$ws.Range("D12").Value = "'0.0000282"
$ws.Range("E12").Value = "  +2.89%  "
$ws.Range("D13").Value = "'679.17"
$ws.Range("E13").Value = "  -1.86%  "
$ws.Range("D14").Value = "3.961.73"
$ws.Range("E14").Value = "  +1.47%  "
$ws.Range("D15").Value = "'8.62"
$ws.Range("E15").Value = "  +1.84%  "
$ws.Range("D16").Value = "69.319.13"
$ws.Range("E16").Value = "  +1.29%  "
$ws.Range("D17").Value = "3.413.28"
$ws.Range("E17").Value = "  +1.84%  "
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("D19").Value = "'17.75"
$ws.Range("E19").Value = "  +1.42%  "
$ws.Range("D20").Value = "'11.32"
$ws.Range("E20").Value = "  +0.70%  "
$ws.Range("D21").Value = "'0.911"
$ws.Range("E21").Value = "  +1.75%  "
$ws.Range("D22").Value = "'5.37"
$ws.Range("E22").Value = "  -2.15%  "
$ws.Range("D23").Value = "'17.03"
$ws.Range("E23").Value = "  +0.45%  "
$ws.Range("D24").Value = "'100.90"
$ws.Range("E24").Value = "  +0.82%  "
$ws.Range("D25").Value = "'3.89"
$ws.Range("E25").Value = "  -0.43%  "
$ws.Range("D26").Value = "'2.69"
$ws.Range("E26").Value = "  -0.60%  "
$ws.Range("D27").Value = "'9.65"
$ws.Range("E27").Value = "  +1.39%  "
$ws.Range("D28").Value = "'33.55"
$ws.Range("E28").Value = "  +1.50%  "
$ws.Range("D29").Value = "'8.75"
$ws.Range("E29").Value = "  +2.36%  "
$ws.Range("D30").Value = "'6.87"
$ws.Range("E30").Value = "  -1.57%  "
$ws.Range("D31").Value = "'3.71"
$ws.Range("E31").Value = "  +8.52%  "
$ws.Range("D32").Value = "'10.99"
$ws.Range("E32").Value = "  -0.93%  "
$ws.Range("D33").Value = "'550.16"
$ws.Range("E33").Value = "  -0.40%  "
$ws.Range("E34").Value = "  -0.17%  "
$ws.Range("D35").Value = "'58.00"
$ws.Range("E35").Value = "  -0.18%  "
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("D37").Value = "3.608.72"
$ws.Range("E37").Value = "  -2.99%  "
$ws.Range("D38").Value = "'0.141"
$ws.Range("E38").Value = "  -0.34%  "
$ws.Range("D39").Value = "'35.03"
$ws.Range("E39").Value = "  +0.86%  "
$ws.Range("D40").Value = "0.0₃0736"
$ws.Range("E40").Value = "  +9.59%  "
$ws.Range("D41").Value = "'3.27"
$ws.Range("E41").Value = "  +2.62%  "
$ws.Range("D42").Value = "'2.69"
$ws.Range("E42").Value = "  +2.66%  "
$ws.Range("D43").Value = "'3.38"
$ws.Range("E43").Value = "  +3.68%  "
$ws.Range("D44").Value = "'0.0424"
$ws.Range("E44").Value = "  +2.68%  "
$ws.Range("D45").Value = "'0.334"
$ws.Range("E45").Value = "  -0.32%  "
$ws.Range("D46").Value = "'2.66"
$ws.Range("E46").Value = "  +0.42%  "
$ws.Range("D47").Value = "'0.129"
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("D48").Value = "'1.40"
$ws.Range("E48").Value = "  +4.80%  "
$ws.Range("E49").Value = "  -0.24%  "
$ws.Range("D50").Value = "'130.87"
$ws.Range("E50").Value = "  -0.76%  "
$ws.Range("D51").Value = "'2.72"
$ws.Range("E51").Value = "  +3.67%  "
